$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the kategori_id value for row 2 (was 11, should be 6)
$ws.Range("A2").Value = 6

# Update the active selection to match the author's final cursor position
$ws.Range("B3").Select()
